# edit.ps1 - Applies the "Add files via upload" commit to docs/Jobs.xlsx
# Adds 24 new job postings (rows 737-760), backfills Date_Applied (I) and
# Rejection/Viewed email columns for a handful of existing rows, converts
# several individual URL formulas into shared-formula groups (matching
# Excel's own fill-down behaviour), and grows Table2 to cover the new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------
# 1. Backfill column I (Rejection_Email date) for rows that got a
#    rejection notice after the original save.
# ---------------------------------------------------------------
$ws.Range("I694").Value = 44210
$ws.Range("I695").Value = 44210
$ws.Range("I696").Value = 44210
$ws.Range("I697").Value = 44210
$ws.Range("I698").Value = 44210
$ws.Range("I712").Value = 44208
$ws.Range("I714").Value = 44208
$ws.Range("I715").Value = 44206

# Column J (Viewed_Email date) for row 732.
$ws.Range("J732").Value = 44207

# ---------------------------------------------------------------
# 2. Re-enter the existing M724:M736 URL formulas as one block so the
#    engine collapses them into a shared formula group, exactly like
#    Excel does when you fill a formula down a column.
# ---------------------------------------------------------------
$ws.Range("M724:M736").Formula = '="https://www.linkedin.com/jobs/search/?currentJobId=" & L724'


# ---------------------------------------------------------------
# 3. Append 24 new job-posting rows (737-760).
# ---------------------------------------------------------------
# ---- New rows 737-760: data cells (A-L), M formulas handled separately below ----
# Row 737
$ws.Range("A737").Value = 'Data Scientist'
$ws.Range("B737").Value = 'Wanted'
$ws.Range("C737").Value = '51-200'
$ws.Range("D737").Value = 'Boston'
$ws.Range("E737").Value = 'MA'
$ws.Range("F737").Formula = "=VLOOKUP(E737,Sheet2!`$A`$1:`$B`$76, 2, FALSE)"
$ws.Range("G737").Value = 44207
$ws.Range("H737").Value = 44207
$ws.Range("J737").Value = 44207
$ws.Range("L737").Value = 2356993969

# Row 738
$ws.Range("A738").Value = 'Entry Level Data Analyst'
$ws.Range("B738").Value = 'KGS Technology Group, Inc'
$ws.Range("C738").Value = '51-200'
$ws.Range("D738").Value = 'West Shokan'
$ws.Range("E738").Value = 'NY'
$ws.Range("F738").Formula = "=VLOOKUP(E738,Sheet2!`$A`$1:`$B`$76, 2, FALSE)"
$ws.Range("G738").Value = 44204
$ws.Range("H738").Value = 44207
$ws.Range("L738").Value = 2372694423

# Row 739
$ws.Range("A739").Value = 'Software Data Engineer'
$ws.Range("B739").Value = 'Hire Talent'
$ws.Range("C739").Value = '501-1000'
$ws.Range("D739").Value = 'Frederick'
$ws.Range("E739").Value = 'MD'
$ws.Range("F739").Formula = "=VLOOKUP(E739,Sheet2!`$A`$1:`$B`$76, 2, FALSE)"
$ws.Range("G739").Value = ">1 week"
$ws.Range("H739").Value = 44207
$ws.Range("L739").Value = 2359239610

# Row 740
$ws.Range("A740").Value = 'Data Analyst'
$ws.Range("B740").Value = 'Pangaea'
$ws.Range("C740").Value = '51-200'
$ws.Range("D740").Value = 'Los Angeles'
$ws.Range("E740").Value = 'CA'
$ws.Range("F740").Formula = "=VLOOKUP(E740,Sheet2!`$A`$1:`$B`$76, 2, FALSE)"
$ws.Range("G740").Value = ">1 week"
$ws.Range("H740").Value = 44207
$ws.Range("L740").Value = 2364729103

# Row 741
$ws.Range("A741").Value = 'Data Analyst - Corporate Housing Technology'
$ws.Range("B741").Value = 'Elliot Scott HR'
$ws.Range("C741").Value = '11-50'
$ws.Range("D741").Value = 'Fort Lauderdale'
$ws.Range("E741").Value = 'FL'
$ws.Range("F741").Formula = "=VLOOKUP(E741,Sheet2!`$A`$1:`$B`$76, 2, FALSE)"
$ws.Range("G741").Value = 44204
$ws.Range("H741").Value = 44207
$ws.Range("J741").Value = 44208
$ws.Range("L741").Value = 2360706212

# Row 742
$ws.Range("A742").Value = 'Data Scientist'
$ws.Range("B742").Value = 'Maxonic'
$ws.Range("C742").Value = '51-200'
$ws.Range("D742").Value = 'Pleasanton'
$ws.Range("E742").Value = 'CA'
$ws.Range("F742").Formula = "=VLOOKUP(E742,Sheet2!`$A`$1:`$B`$76, 2, FALSE)"
$ws.Range("G742").Value = 44202
$ws.Range("H742").Value = 44207
$ws.Range("J742").Value = 44208
$ws.Range("L742").Value = 2355287426

# Row 743
$ws.Range("A743").Value = 'Data Analyst'
$ws.Range("B743").Value = 'Tucker Parker Smith Group'
$ws.Range("C743").Value = '11-50'
$ws.Range("D743").Value = 'San Francisco'
$ws.Range("E743").Value = 'CA'
$ws.Range("F743").Formula = "=VLOOKUP(E743,Sheet2!`$A`$1:`$B`$76, 2, FALSE)"
$ws.Range("G743").Value = 44206
$ws.Range("H743").Value = 44207
$ws.Range("L743").Value = 2362798414

# Row 744
$ws.Range("A744").Value = 'Data Analyst I '
$ws.Range("B744").Value = 'Hire Talent'
$ws.Range("C744").Value = '501-1000'
$ws.Range("D744").Value = 'San Francisco'
$ws.Range("E744").Value = 'CA'
$ws.Range("F744").Formula = "=VLOOKUP(E744,Sheet2!`$A`$1:`$B`$76, 2, FALSE)"
$ws.Range("G744").Value = 44204
$ws.Range("H744").Value = 44207
$ws.Range("L744").Value = 2372101133

# Row 745
$ws.Range("A745").Value = 'Data Analyst'
$ws.Range("B745").Value = 'US Tech Solutions'
$ws.Range("C745").Value = '1001-5000'
$ws.Range("D745").Value = 'Washington'
$ws.Range("E745").Value = 'DC'
$ws.Range("F745").Formula = "=VLOOKUP(E745,Sheet2!`$A`$1:`$B`$76, 2, FALSE)"
$ws.Range("G745").Value = 44204
$ws.Range("H745").Value = 44207
$ws.Range("L745").Value = 2372193668

# Row 746
$ws.Range("A746").Value = 'Python Developer/ML/Data Engineer'
$ws.Range("B746").Value = 'X-Team'
$ws.Range("D746").Value = 'Remote'
$ws.Range("G746").Value = ">1 week"
$ws.Range("H746").Value = 44207
$ws.Range("J746").Value = 44208

# Row 747
$ws.Range("A747").Value = 'Data Engineering'
$ws.Range("B747").Value = 'Hire Talent'
$ws.Range("C747").Value = '501-1000'
$ws.Range("D747").Value = 'Beaverton'
$ws.Range("E747").Value = 'OR'
$ws.Range("F747").Formula = "=VLOOKUP(E747,Sheet2!`$A`$1:`$B`$76, 2, FALSE)"
$ws.Range("G747").Value = 44207
$ws.Range("H747").Value = 44209
$ws.Range("L747").Value = 2376748343

# Row 748
$ws.Range("A748").Value = 'Data Analyst'
$ws.Range("B748").Value = 'Good Apple'
$ws.Range("C748").Value = '51-200'
$ws.Range("D748").Value = 'New York'
$ws.Range("E748").Value = 'NY'
$ws.Range("F748").Formula = "=VLOOKUP(E748,Sheet2!`$A`$1:`$B`$76, 2, FALSE)"
$ws.Range("G748").Value = 44208
$ws.Range("H748").Value = 44209
$ws.Range("L748").Value = 2359609615

# Row 749
$ws.Range("A749").Value = 'Data Engineer'
$ws.Range("B749").Value = 'Burtch Works'
$ws.Range("C749").Value = '11-50'
$ws.Range("D749").Value = 'Boston'
$ws.Range("E749").Value = 'MA'
$ws.Range("F749").Formula = "=VLOOKUP(E749,Sheet2!`$A`$1:`$B`$76, 2, FALSE)"
$ws.Range("G749").Value = 44207
$ws.Range("H749").Value = 44209
$ws.Range("L749").Value = 2350296378

# Row 750
$ws.Range("A750").Value = 'Data Analyst (risk)'
$ws.Range("B750").Value = 'Harnham'
$ws.Range("C750").Value = '51-200'
$ws.Range("D750").Value = 'San Francisco'
$ws.Range("E750").Value = 'CA'
$ws.Range("F750").Formula = "=VLOOKUP(E750,Sheet2!`$A`$1:`$B`$76, 2, FALSE)"
$ws.Range("G750").Value = 44208
$ws.Range("H750").Value = 44209
$ws.Range("L750").Value = 2368083278

# Row 751
$ws.Range("A751").Value = 'Data Analyst'
$ws.Range("B751").Value = 'Apex Systems'
$ws.Range("C751").Value = '1001-5000'
$ws.Range("D751").Value = 'Mountain View'
$ws.Range("E751").Value = 'CA'
$ws.Range("F751").Formula = "=VLOOKUP(E751,Sheet2!`$A`$1:`$B`$76, 2, FALSE)"
$ws.Range("G751").Value = 44209
$ws.Range("H751").Value = 44209
$ws.Range("L751").Value = 2366585560

# Row 752
$ws.Range("A752").Value = 'Product Data Analyst- Analyze the "Healthiverse"'
$ws.Range("B752").Value = 'DrFirst, Inc.'
$ws.Range("C752").Value = '201-500'
$ws.Range("D752").Value = 'Rockville'
$ws.Range("E752").Value = 'MD'
$ws.Range("F752").Formula = "=VLOOKUP(E752,Sheet2!`$A`$1:`$B`$76, 2, FALSE)"
$ws.Range("G752").Value = 44208
$ws.Range("H752").Value = 44209
$ws.Range("L752").Value = 2366505315

# Row 753
$ws.Range("A753").Value = 'Data Analyst'
$ws.Range("B753").Value = 'Brooksource'
$ws.Range("C753").Value = '1001-5000'
$ws.Range("D753").Value = 'New York'
$ws.Range("E753").Value = 'NY'
$ws.Range("F753").Formula = "=VLOOKUP(E753,Sheet2!`$A`$1:`$B`$76, 2, FALSE)"
$ws.Range("G753").Value = 44209
$ws.Range("H753").Value = 44210
$ws.Range("L753").Value = 2377936621

# Row 754
$ws.Range("A754").Value = 'Data Analyst II (IT)'
$ws.Range("B754").Value = 'Hire Talent'
$ws.Range("C754").Value = '501-1000'
$ws.Range("D754").Value = 'Framingham'
$ws.Range("E754").Value = 'MA'
$ws.Range("F754").Formula = "=VLOOKUP(E754,Sheet2!`$A`$1:`$B`$76, 2, FALSE)"
$ws.Range("G754").Value = 44205
$ws.Range("H754").Value = 44210
$ws.Range("L754").Value = 2373270367

# Row 755
$ws.Range("A755").Value = 'Data Analyst, Analytics (Contractor)'
$ws.Range("B755").Value = 'Evernote'
$ws.Range("C755").Value = '201-500'
$ws.Range("D755").Value = 'San Diego'
$ws.Range("E755").Value = 'CA'
$ws.Range("F755").Formula = "=VLOOKUP(E755,Sheet2!`$A`$1:`$B`$76, 2, FALSE)"
$ws.Range("G755").Value = 44209
$ws.Range("H755").Value = 44210
$ws.Range("L755").Value = 2369175427

# Row 756
$ws.Range("A756").Value = 'Data Analyst (Merchandising)'
$ws.Range("B756").Value = 'FIGS'
$ws.Range("C756").Value = '51-200'
$ws.Range("D756").Value = 'Santa Monica'
$ws.Range("E756").Value = 'CA'
$ws.Range("F756").Formula = "=VLOOKUP(E756,Sheet2!`$A`$1:`$B`$76, 2, FALSE)"
$ws.Range("G756").Value = 44210
$ws.Range("H756").Value = 44210
$ws.Range("L756").Value = 2369186406

# Row 757
$ws.Range("A757").Value = 'Data Scientist I'
$ws.Range("B757").Value = 'Caterpillar'
$ws.Range("C757").Value = '10001'
$ws.Range("D757").Value = 'Peoria'
$ws.Range("E757").Value = 'IL'
$ws.Range("F757").Formula = "=VLOOKUP(E757,Sheet2!`$A`$1:`$B`$76, 2, FALSE)"
$ws.Range("G757").Value = 44210
$ws.Range("H757").Value = 44210
$ws.Range("L757").Value = 2368043859

# Row 758
$ws.Range("A758").Value = 'Associate Data Scientist'
$ws.Range("B758").Value = 'The Home Depot'
$ws.Range("C758").Value = '10001'
$ws.Range("D758").Value = 'Atlanta'
$ws.Range("E758").Value = 'GA'
$ws.Range("F758").Formula = "=VLOOKUP(E758,Sheet2!`$A`$1:`$B`$76, 2, FALSE)"
$ws.Range("G758").Value = 44209
$ws.Range("H758").Value = 44210
$ws.Range("L758").Value = 2359625710

# Row 759
$ws.Range("A759").Value = 'Data Reporting Analyst'
$ws.Range("B759").Value = 'SNI Financial'
$ws.Range("C759").Value = '201-500'
$ws.Range("D759").Value = 'Tampa'
$ws.Range("E759").Value = 'FL'
$ws.Range("F759").Formula = "=VLOOKUP(E759,Sheet2!`$A`$1:`$B`$76, 2, FALSE)"
$ws.Range("G759").Value = 44210
$ws.Range("H759").Value = 44210
$ws.Range("L759").Value = 2380952280

# Row 760
$ws.Range("A760").Value = 'Data Analyst'
$ws.Range("B760").Value = 'EPITEC'
$ws.Range("C760").Value = '1001-5000'
$ws.Range("D760").Value = 'Grand Blanc'
$ws.Range("E760").Value = 'MI'
$ws.Range("F760").Formula = "=VLOOKUP(E760,Sheet2!`$A`$1:`$B`$76, 2, FALSE)"
$ws.Range("G760").Value = 44210
$ws.Range("H760").Value = 44210
$ws.Range("L760").Value = 2369603611


# ---------------------------------------------------------------
# 4. URL formulas (column M) for the new rows, written in the same
#    three blocks the author's own fill-down produced: a lone cell
#    for row 737 (copied from a single source cell, not a fill), then
#    three shared-formula groups for 738-745, 746-752, 753-760.
# ---------------------------------------------------------------
$ws.Range("M737").Formula = '="https://www.linkedin.com/jobs/search/?currentJobId=" & L737'
$ws.Range("M738:M745").Formula = '="https://www.linkedin.com/jobs/search/?currentJobId=" & L738'
$ws.Range("M746:M752").Formula = '="https://www.linkedin.com/jobs/search/?currentJobId=" & L746'
$ws.Range("M753:M760").Formula = '="https://www.linkedin.com/jobs/search/?currentJobId=" & L753'

# ---------------------------------------------------------------
# 5. Grow Table2 (and its AutoFilter) to cover the new rows.
# ---------------------------------------------------------------
$tbl = $ws.ListObjects.Item("Table2")
$tbl.Resize($ws.Range("A1:M760"))

# ---------------------------------------------------------------
# 6. Leave the selection where the author ended up (freeze pane stays
#    at row 1; the visible window was scrolled down to row 728 and the
#    active cell left on the newly-added D741).
# ---------------------------------------------------------------
$ws.Range("D741").Select()
$excel.ActiveWindow.ScrollRow = 728
